$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Week-of-Jan-4 section (rows 59-68): correct an hours entry and log a
# previously-blank activity row
$ws.Range("C61").Value = 2.5

$ws.Range("A62").Value = "High-level review of anomaly detection algorithms"
$ws.Range("B62").Value = "Background Review"
$ws.Range("C62").Value = 1.5

# Week-of-Jan-11 section (rows 70-79): log two previously-blank activity rows
$ws.Range("A70").Value = "Finalize proposal and presentation"
$ws.Range("B70").Value = "Project Documents"
$ws.Range("C70").Value = 1

$ws.Range("A71").Value = "Update GitHub READMEs"
$ws.Range("B71").Value = "Project Management"
$ws.Range("C71").Value = 0.5

# Leave the sheet scrolled/selected where the user last left off
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("H96").Select()
